$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "42.23") are stored as text, not converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.584.71'
$ws.Range("E2").Value = '  +0.88%  '
$ws.Range("D3").Value = '1.872.87'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -1.09%  '
$ws.Range("D5").Value = '314.36'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '0.5076'
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("D8").Value = '0.3909'
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("D9").Value = '0.08355'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").Value = '42.23'
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = '1.106'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '6.197'
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '1.871.70'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").Value = '20.32'
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = '7.259'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").Value = '93.18'
$ws.Range("E17").Value = '  +3.03%  '
$ws.Range("D18").Value = '0.00001099'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '0.06717'
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '17.62'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("D22").Value = '5.929'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").Value = '28.608.04'
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").Value = '11.07'
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").Value = '2.190'
$ws.Range("E25").Value = '  -4.03%  '
$ws.Range("D26").Value = '2.081.60'
$ws.Range("E26").Value = '  +2.06%  '
$ws.Range("D27").Value = '157.85'
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("D28").Value = '20.56'
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("D29").Value = '2.423'
$ws.Range("E29").Value = '  +3.02%  '
$ws.Range("D30").Value = '126.57'
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = '0.1037'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").Value = '1.047'
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("D33").Value = '5.775'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '3.633'
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").Value = '0.02443'
$ws.Range("E35").Value = '  +1.02%  '
$ws.Range("D36").Value = '0.06571'
$ws.Range("E36").Value = '  +1.75%  '
$ws.Range("D37").Value = '9.017'
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("D38").Value = '0.2163'
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("D39").Value = '5.043'
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").Value = '1.188'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '1.239'
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("D42").Value = '0.6373'
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("D43").Value = '11.13'
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D44").Value = '1.005'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("D45").Value = '0.5985'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("D46").Value = '13.00'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").Value = '3.675'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '2.006'
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").Value = '1.210'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").Value = '122.27'
$ws.Range("E50").Value = '  +1.21%  '
$ws.Range("D51").Value = '1.176'
$ws.Range("E51").Value = '  -2.88%  '

# Restore default cell style for column D now that the text values are set,
# so no stray number-format style is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
